# Update prev_task_id / next_task_id columns (I, J) so that the operation
# name segment (2nd underscore-delimited field) of each task id is replaced
# with its corresponding numeric operation code.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "PPF점착"      = "23311"
    "TOP COATING"  = "12400"
    "인쇄"          = "20906"
    "염료점착"      = "20300"
    "안료접착"      = "20902"
    "투명점착"      = "20500"
    "유광 S/R"      = "20700"
    "MIBK SR"      = "20706"
}

$cols = @("I", "J")
for ($row = 2; $row -le 14; $row++) {
    foreach ($col in $cols) {
        $addr = "$col$row"
        $val = $ws.Range($addr).Value()
        if ($null -eq $val) { continue }

        $parts = $val -split "_"
        if ($parts.Length -ge 2) {
            $opName = $parts[1]
            if ($map.ContainsKey($opName)) {
                $parts[1] = $map[$opName]
                $newVal = [string]::Join("_", $parts)
                $ws.Range($addr).Value = $newVal
            }
        }
    }
}
